# "Delete redundant old signs"
#
# The deck of signs is re-ordered: two new signs ("ARENA" and "HELPDESK"),
# built the same way as the existing "TOILETS" sign (slide 1 - a duplicate
# of that slide with only the text changed), are inserted right after the
# "TOILETS" slide. The old, now-redundant "ARENA" sign slide is removed.

$p = $ppt.ActivePresentation

# --- Build the two new sign slides by duplicating the "TOILETS" slide ---
# (slide 1) which already uses the correct layout/placeholder for a
# text-only sign. Duplicate() inserts the copy immediately after its
# source, so duplicating the freshly-made "ARENA" slide places "HELPDESK"
# right after it, keeping TOILETS, ARENA, HELPDESK in order.

$toilets = $p.Slides.Item(1)

$arenaRange = $toilets.Duplicate()
$arenaSlide = $arenaRange.Item(1)
$arenaSlide.Shapes.Item(1).TextFrame.TextRange.Text = "ARENA"

$helpdeskRange = $arenaSlide.Duplicate()
$helpdeskSlide = $helpdeskRange.Item(1)
$helpdeskSlide.Shapes.Item(1).TextFrame.TextRange.Text = "HELPDESK"

# --- Remove the old, now-redundant "ARENA" sign slide ---
# After the inserts above the slide order is:
#   1 TOILETS, 2 ARENA(new), 3 HELPDESK(new), 4 PHOTOGRAPHY ..., 5 EXIT,
#   6 ARENA(old), 7 STAGING IN, 8 STAGING OUT, 9 QUIET ZONE,
#   10 POWER TOOLS AREA, 11 NO ENTRY
# so the old "ARENA" sign is now at index 6. Find it by its text rather
# than hard-coding the index, in case anything above shifts.

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.Shapes.Count -eq 1 -and $slide.Shapes.Item(1).HasTextFrame) {
        $text = $slide.Shapes.Item(1).TextFrame.TextRange.Text
        if ($text -eq "ARENA" -and $slide.SlideIndex -ne $arenaSlide.SlideIndex) {
            $slide.Delete()
            break
        }
    }
}
